$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 10.04013707806289
$ws.Cells.Item(2, 3).Value = 4.899733097183971
$ws.Cells.Item(2, 5).Value = 22.90544527475214
$ws.Cells.Item(2, 6).Value = 41.17848498620837
$ws.Cells.Item(2, 7).Value = 3.640622344164938
$ws.Cells.Item(2, 10).Value = 8.133713967735492
$ws.Cells.Item(2, 11).Value = 9.337930198256901
$ws.Cells.Item(2, 14).Value = 18.38961945010614
$ws.Cells.Item(2, 15).Value = 22.49223145042268

$ws.Cells.Item(3, 2).Value = 9.763860862809988
$ws.Cells.Item(3, 3).Value = 4.735599477558031
$ws.Cells.Item(3, 5).Value = 22.58342008695665
$ws.Cells.Item(3, 6).Value = 40.99659910725056
$ws.Cells.Item(3, 7).Value = 3.64230783541312
$ws.Cells.Item(3, 10).Value = 8.156460165211362
$ws.Cells.Item(3, 11).Value = 9.149467944207627
$ws.Cells.Item(3, 14).Value = 18.4486136099298
$ws.Cells.Item(3, 15).Value = 22.58257689046113

$ws.Cells.Item(4, 2).Value = 9.591952763366249
$ws.Cells.Item(4, 3).Value = 4.6310937063305
$ws.Cells.Item(4, 5).Value = 22.38902302831827
$ws.Cells.Item(4, 6).Value = 40.89597975989211
$ws.Cells.Item(4, 7).Value = 3.643396811284474
$ws.Cells.Item(4, 10).Value = 8.171266979730721
$ws.Cells.Item(4, 11).Value = 9.033201312167343
$ws.Cells.Item(4, 14).Value = 18.48658814554508
$ws.Cells.Item(4, 15).Value = 22.64284272057421

$ws.Cells.Item(5, 2).Value = 9.521436653549014
$ws.Cells.Item(5, 3).Value = 4.587610891009387
$ws.Cells.Item(5, 5).Value = 22.31073739132879
$ws.Cells.Item(5, 6).Value = 40.85778761579161
$ws.Cells.Item(5, 7).Value = 3.643854219113667
$ws.Cells.Item(5, 10).Value = 8.177512672188167
$ws.Cells.Item(5, 11).Value = 8.985747749734212
$ws.Cells.Item(5, 14).Value = 18.50250487700847
$ws.Cells.Item(5, 15).Value = 22.668604437214

$ws.Cells.Item(6, 2).Value = 9.509702914066024
$ws.Cells.Item(6, 3).Value = 4.580337787704339
$ws.Cells.Item(6, 5).Value = 22.2977972273742
$ws.Cells.Item(6, 6).Value = 40.85161641799696
$ws.Cells.Item(6, 7).Value = 3.643930996505582
$ws.Cells.Item(6, 10).Value = 8.178562569300365
$ws.Cells.Item(6, 11).Value = 8.977865567156261
$ws.Cells.Item(6, 14).Value = 18.50517455933235
$ws.Cells.Item(6, 15).Value = 22.67295473185894

$ws.Cells.Item(7, 2).Value = 9.59100348052327
$ws.Cells.Item(7, 3).Value = 4.630510851396123
$ws.Cells.Item(7, 5).Value = 22.38796333868552
$ws.Cells.Item(7, 6).Value = 40.89545326883828
$ws.Cells.Item(7, 7).Value = 3.643402924755704
$ws.Cells.Item(7, 10).Value = 8.171350353197843
$ws.Cells.Item(7, 11).Value = 9.032561550088449
$ws.Cells.Item(7, 14).Value = 18.48680101363576
$ws.Cells.Item(7, 15).Value = 22.6431852848312

$ws.Cells.Item(8, 2).Value = 9.945412622426238
$ws.Cells.Item(8, 3).Value = 4.843936848807921
$ws.Cells.Item(8, 5).Value = 22.79377835391276
$ws.Cells.Item(8, 6).Value = 41.11349468493228
$ws.Cells.Item(8, 7).Value = 3.641192303781259
$ws.Cells.Item(8, 10).Value = 8.141382657795114
$ws.Cells.Item(8, 11).Value = 9.273100634277238
$ws.Cells.Item(8, 14).Value = 18.40959780257157
$ws.Cells.Item(8, 15).Value = 22.52238623898458

$ws.Cells.Item(9, 2).Value = 10.61788843526886
$ws.Cells.Item(9, 3).Value = 5.231290878589669
$ws.Cells.Item(9, 5).Value = 23.61170934728518
$ws.Cells.Item(9, 6).Value = 41.62730994411086
$ws.Cells.Item(9, 7).Value = 3.637284391108655
$ws.Cells.Item(9, 10).Value = 8.089266034455502
$ws.Cells.Item(9, 11).Value = 9.737703102741431
$ws.Cells.Item(9, 14).Value = 18.27204411967526
$ws.Cells.Item(9, 15).Value = 22.32364118245237

$ws.Cells.Item(10, 2).Value = 11.09271853324387
$ws.Cells.Item(10, 3).Value = 5.495012919508892
$ws.Cells.Item(10, 5).Value = 24.22038687174014
$ws.Cells.Item(10, 6).Value = 42.05509015253107
$ws.Cells.Item(10, 7).Value = 3.634670857609474
$ws.Cells.Item(10, 10).Value = 8.055002814963515
$ws.Cells.Item(10, 11).Value = 10.0712333986749
$ws.Cells.Item(10, 14).Value = 18.17933834480635
$ws.Cells.Item(10, 15).Value = 22.20102302258388

$ws.Cells.Item(11, 2).Value = 11.303533413796
$ws.Cells.Item(11, 3).Value = 5.610135367457848
$ws.Cells.Item(11, 5).Value = 24.49776608140973
$ws.Cells.Item(11, 6).Value = 42.26005171249606
$ws.Cells.Item(11, 7).Value = 3.633537242035925
$ws.Cells.Item(11, 10).Value = 8.040284297860442
$ws.Cells.Item(11, 11).Value = 10.22058550491946
$ws.Cells.Item(11, 14).Value = 18.13896058807044
$ws.Cells.Item(11, 15).Value = 22.15035352953944

$ws.Cells.Item(12, 2).Value = 11.3825423333211
$ws.Cells.Item(12, 3).Value = 5.65300993015058
$ws.Cells.Item(12, 5).Value = 24.60277433518296
$ws.Cells.Item(12, 6).Value = 42.3391028046981
$ws.Cells.Item(12, 7).Value = 3.633115877969343
$ws.Cells.Item(12, 10).Value = 8.034835176281399
$ws.Cells.Item(12, 11).Value = 10.27674739216792
$ws.Cells.Item(12, 14).Value = 18.12392733494878
$ws.Cells.Item(12, 15).Value = 22.13190379531052

$ws.Cells.Item(13, 2).Value = 11.36556404864143
$ws.Cells.Item(13, 3).Value = 5.643808474141351
$ws.Cells.Item(13, 5).Value = 24.5801617849338
$ws.Cells.Item(13, 6).Value = 42.32201472667158
$ws.Cells.Item(13, 7).Value = 3.633206275020091
$ws.Cells.Item(13, 10).Value = 8.036003212237254
$ws.Cells.Item(13, 11).Value = 10.26467028486185
$ws.Cells.Item(13, 14).Value = 18.12715360591833
$ws.Cells.Item(13, 15).Value = 22.13584442058374

$ws.Cells.Item(14, 2).Value = 11.31005041205264
$ws.Cells.Item(14, 3).Value = 5.613677230668901
$ws.Cells.Item(14, 5).Value = 24.50640632002425
$ws.Cells.Item(14, 6).Value = 42.26652681847316
$ws.Cells.Item(14, 7).Value = 3.633502417827168
$ws.Cells.Item(14, 10).Value = 8.039833502965912
$ws.Cells.Item(14, 11).Value = 10.22521417752505
$ws.Cells.Item(14, 14).Value = 18.13771865050685
$ws.Cells.Item(14, 15).Value = 22.14882085774154

$ws.Cells.Item(15, 2).Value = 11.2759375071612
$ws.Cells.Item(15, 3).Value = 5.595126602423923
$ws.Cells.Item(15, 5).Value = 24.46122225290076
$ws.Cells.Item(15, 6).Value = 42.23272432316526
$ws.Cells.Item(15, 7).Value = 3.633684842993081
$ws.Cells.Item(15, 10).Value = 8.04219586476402
$ws.Cells.Item(15, 11).Value = 10.20099326881926
$ws.Cells.Item(15, 14).Value = 18.14422347158595
$ws.Cells.Item(15, 15).Value = 22.1568654558661

$ws.Cells.Item(16, 2).Value = 11.07883022584734
$ws.Cells.Item(16, 3).Value = 5.48738985194273
$ws.Cells.Item(16, 5).Value = 24.20226075738546
$ws.Cells.Item(16, 6).Value = 42.04189944300222
$ws.Cells.Item(16, 7).Value = 3.634746051312441
$ws.Cells.Item(16, 10).Value = 8.055982140418974
$ws.Cells.Item(16, 11).Value = 10.06142050983293
$ws.Cells.Item(16, 14).Value = 18.18201314055429
$ws.Cells.Item(16, 15).Value = 22.20443746792476

$ws.Cells.Item(17, 2).Value = 10.95652617060905
$ws.Cells.Item(17, 3).Value = 5.420038373130718
$ws.Cells.Item(17, 5).Value = 24.04344984687016
$ws.Cells.Item(17, 6).Value = 41.9274514729055
$ws.Cells.Item(17, 7).Value = 3.635411202035618
$ws.Cells.Item(17, 10).Value = 8.064661633215312
$ws.Cells.Item(17, 11).Value = 9.975150475321955
$ws.Cells.Item(17, 14).Value = 18.2056547017208
$ws.Cells.Item(17, 15).Value = 22.23493223193619

$ws.Cells.Item(18, 2).Value = 10.88569591087986
$ws.Cells.Item(18, 3).Value = 5.380844625950652
$ws.Cells.Item(18, 5).Value = 23.95215782492309
$ws.Cells.Item(18, 6).Value = 41.8626027655493
$ws.Cells.Item(18, 7).Value = 3.635798986015397
$ws.Cells.Item(18, 10).Value = 8.06973556834196
$ws.Cells.Item(18, 11).Value = 9.925309904923015
$ws.Cells.Item(18, 14).Value = 18.21942168195365
$ws.Cells.Item(18, 15).Value = 22.25295283169483

$ws.Cells.Item(19, 2).Value = 10.86163327598192
$ws.Cells.Item(19, 3).Value = 5.367496905579983
$ws.Cells.Item(19, 5).Value = 23.92125988775834
$ws.Cells.Item(19, 6).Value = 41.84081578724827
$ws.Cells.Item(19, 7).Value = 3.635931178507505
$ws.Cells.Item(19, 10).Value = 8.071467562084228
$ws.Cells.Item(19, 11).Value = 9.908398582706756
$ws.Cells.Item(19, 14).Value = 18.2241120019667
$ws.Cells.Item(19, 15).Value = 22.25913679028056

$ws.Cells.Item(20, 2).Value = 10.96959633308214
$ws.Cells.Item(20, 3).Value = 5.427255315155139
$ws.Cells.Item(20, 5).Value = 24.06035090601105
$ws.Cells.Item(20, 6).Value = 41.93953373690593
$ws.Cells.Item(20, 7).Value = 3.635339857001736
$ws.Cells.Item(20, 10).Value = 8.063729231223371
$ws.Cells.Item(20, 11).Value = 9.984357267963397
$ws.Cells.Item(20, 14).Value = 18.20312053663803
$ws.Cells.Item(20, 15).Value = 22.23163622848065

$ws.Cells.Item(21, 2).Value = 11.32637899580765
$ws.Cells.Item(21, 3).Value = 5.622547211119524
$ws.Cells.Item(21, 5).Value = 24.52807166112881
$ws.Cells.Item(21, 6).Value = 42.28278641384172
$ws.Cells.Item(21, 7).Value = 3.633415219103366
$ws.Cells.Item(21, 10).Value = 8.038705078679959
$ws.Cells.Item(21, 11).Value = 10.23681450826325
$ws.Cells.Item(21, 14).Value = 18.13460847681377
$ws.Cells.Item(21, 15).Value = 22.1449893197366

$ws.Cells.Item(22, 2).Value = 11.55473489916045
$ws.Cells.Item(22, 3).Value = 5.745977516942673
$ws.Cells.Item(22, 5).Value = 24.83353965447197
$ws.Cells.Item(22, 6).Value = 42.51546868992169
$ws.Cells.Item(22, 7).Value = 3.632203453354985
$ws.Cells.Item(22, 10).Value = 8.023075644806465
$ws.Cells.Item(22, 11).Value = 10.39948800694405
$ws.Cells.Item(22, 14).Value = 18.09132893856737
$ws.Cells.Item(22, 15).Value = 22.09266177319312

$ws.Cells.Item(23, 2).Value = 11.4333216133752
$ws.Cells.Item(23, 3).Value = 5.680491966308565
$ws.Cells.Item(23, 5).Value = 24.6705574313572
$ws.Cells.Item(23, 6).Value = 42.39053644628014
$ws.Cells.Item(23, 7).Value = 3.632845990741219
$ws.Cells.Item(23, 10).Value = 8.031351118683975
$ws.Cells.Item(23, 11).Value = 10.3128953112023
$ws.Cells.Item(23, 14).Value = 18.11429143163932
$ws.Cells.Item(23, 15).Value = 22.12019546370474

$ws.Cells.Item(24, 2).Value = 10.96368891406691
$ws.Cells.Item(24, 3).Value = 5.423994005615303
$ws.Cells.Item(24, 5).Value = 24.05270989686706
$ws.Cells.Item(24, 6).Value = 41.93406838529158
$ws.Cells.Item(24, 7).Value = 3.635372095329584
$ws.Cells.Item(24, 10).Value = 8.064150508547522
$ws.Cells.Item(24, 11).Value = 9.980195629097805
$ws.Cells.Item(24, 14).Value = 18.20426568703999
$ws.Cells.Item(24, 15).Value = 22.23312482938026

$ws.Cells.Item(25, 2).Value = 10.43897364014803
$ws.Cells.Item(25, 3).Value = 5.130045636430437
$ws.Cells.Item(25, 5).Value = 23.38866325819232
$ws.Cells.Item(25, 6).Value = 41.47929874716537
$ws.Cells.Item(25, 7).Value = 3.638296147268955
$ws.Cells.Item(25, 10).Value = 8.089266034455502
$ws.Cells.Item(25, 11).Value = 9.737703102741431
$ws.Cells.Item(25, 14).Value = 18.30778286002454
$ws.Cells.Item(25, 15).Value = 22.37330839022836
